# CPI corrected, new Plots
# Update the CSCC column (C) values for rows 2-5 with corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 304.3919333712181
$ws.Range("C3").Value = 69.62955639839001
$ws.Range("C4").Value = 41.31740201574306
$ws.Range("C5").Value = 264.1145392529322
